# McCarter & Ash Reading Test 2 Practiced.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Clear the placeholder "Speaking" (K) score of 5 for rows that have not
#     actually been practiced yet. This recalculates the dependent "Overall"
#     (L) average formulas automatically.
$kRowsToClear = @(4, 5, 8, 9, 10, 11, 12, 13, 26)
foreach ($r in $kRowsToClear) {
    $ws.Range("K$r").Value = ""
}

# --- Add the new test entry: McCarter&Ash Test 2 (row 38), Reading practiced.
# Pull formatting for the new cells from the row above (row 37), which is the
# most recently-filled data row, so number formats / fonts / borders match.
$ws.Range("C37").Copy()
$ws.Range("C38").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D37").Copy()
$ws.Range("D38").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I37").Copy()
$ws.Range("I38").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C38").Value = 35
$ws.Range("D38").Value = 45520
$ws.Range("E38").Value = "McCarter&Ash Test 2"
$ws.Range("H38").Value = 25
$ws.Range("I38").Formula = "=IFERROR(INDEX(Sheet2!`$F`$5:`$F`$20, MATCH(Table1[[#This Row],[Read_Mark]], Sheet2!`$D`$5:`$D`$20, 1)),""No Grade"")"

# --- Update the sheet view: scroll position and current selection.
[void]$ws.Range("D39").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
